$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Append three new rows (104-106) describing new SingleR RData temp folders
# produced while finishing the A/C/N pipe_23-06 run.
$ws.Range("A104").Value = "Kriegstein"
$ws.Range("B104").Value = "SingleR_RData_2022-07-04 11-47-53"
$ws.Range("C104").Value = "Kriegstein to SingleR"
$ws.Range("D104").Value = "A"

$ws.Range("A105").Value = "Kriegstein"
$ws.Range("B105").Value = "SingleR_RData_2022-07-04 11-49-06"
$ws.Range("C105").Value = "Kriegstein to SingleR"
$ws.Range("D105").Value = "C"

$ws.Range("A106").Value = "Kriegstein"
$ws.Range("B106").Value = "SingleR_RData_2022-07-04 11-49-39"
$ws.Range("C106").Value = "Kriegstein to SingleR"
$ws.Range("D106").Value = "N"

# Update the view so the newly added rows are visible, matching the
# author's scroll position/selection after editing the log.
$ws.Range("A107").Select()
$excel.ActiveWindow.ScrollRow = 91
$excel.ActiveWindow.ScrollColumn = 1
